# Delete the data row that contains Record_ID_list_1 = 734 and
# Record_ID_list_2 = 655 (worksheet row 621). Deleting the entire row
# shifts all the rows below it up by one, which shrinks the sheet
# dimension to A1:C629 and drops the last (now-empty) row 630.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the frozen-pane scroll/selection state before editing so the
# final view matches a fresh selection at the top of the sheet.
$ws.Rows.Item(621).Delete()

# The autofilter range (and the corresponding _FilterDatabase defined
# name) still reference the old, now too-large range after the row
# delete; re-apply the autofilter over the shrunk used range so it
# tracks the new extent of the data (A1:C629).
$ws.AutoFilterMode = $false
$ws.Range("A1:C629").AutoFilter() | Out-Null

# Excel does not always refresh the hidden _FilterDatabase name's
# reference when a row is deleted; update it explicitly to match the
# new autofilter extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$629"
    }
}

# Reset the view to the top-left of the data (cell A2), which is what
# the saved workbook reflects after the edit.
$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
